# Updates recomputed transition-probability values in the
# "St. Bonaventure_A" team-specific matrix (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1550387596899225
$ws.Range("C2").Value = 0.6589147286821705
$ws.Range("J2").Value = 0.003875968992248062
$ws.Range("P2").Value = 0.1124031007751938
$ws.Range("S2").Value = 0.06976744186046512
$ws.Range("C3").Value = 0.02840909090909091
$ws.Range("J3").Value = 0.02272727272727273
$ws.Range("P3").Value = 0.7556818181818182
$ws.Range("S3").Value = 0.1931818181818182
$ws.Range("P4").Value = 0.6976744186046512
$ws.Range("S4").Value = 0.3023255813953488
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05
$ws.Range("D6").Value = 0.005
$ws.Range("E6").Value = 0.005
$ws.Range("F6").Value = 0.075
$ws.Range("J6").Value = 0.245
$ws.Range("O6").Value = 0.025
$ws.Range("Q6").Value = 0.14
$ws.Range("R6").Value = 0.075
$ws.Range("S6").Value = 0.38
$ws.Range("B7").Value = 0.12
$ws.Range("D7").Value = 0.01142857142857143
$ws.Range("F7").Value = 0.05714285714285714
$ws.Range("J7").Value = 0.1542857142857143
$ws.Range("O7").Value = 0.01142857142857143
$ws.Range("Q7").Value = 0.2342857142857143
$ws.Range("R7").Value = 0.07428571428571429
$ws.Range("S7").Value = 0.3371428571428571
$ws.Range("B8").Value = 0.116751269035533
$ws.Range("D8").Value = 0.02538071065989848
$ws.Range("F8").Value = 0.0583756345177665
$ws.Range("J8").Value = 0.1192893401015228
$ws.Range("O8").Value = 0.01269035532994924
$ws.Range("Q8").Value = 0.1852791878172589
$ws.Range("R8").Value = 0.09898477157360407
$ws.Range("S8").Value = 0.383248730964467
$ws.Range("B9").Value = 0.0797872340425532
$ws.Range("D9").Value = 0.02127659574468085
$ws.Range("E9").Value = 0.005319148936170213
$ws.Range("F9").Value = 0.05851063829787234
$ws.Range("J9").Value = 0.1117021276595745
$ws.Range("O9").Value = 0.02659574468085106
$ws.Range("Q9").Value = 0.1542553191489362
$ws.Range("R9").Value = 0.0851063829787234
$ws.Range("S9").Value = 0.4574468085106383
$ws.Range("B10").Value = 0.1024590163934426
$ws.Range("D10").Value = 0.02377049180327869
$ws.Range("E10").Value = 0.000819672131147541
$ws.Range("F10").Value = 0.06885245901639345
$ws.Range("J10").Value = 0.1139344262295082
$ws.Range("O10").Value = 0.01065573770491803
$ws.Range("Q10").Value = 0.2270491803278689
$ws.Range("R10").Value = 0.1024590163934426
$ws.Range("S10").Value = 0.35
$ws.Range("G11").Value = 0.1262798634812287
$ws.Range("J11").Value = 0.10580204778157
$ws.Range("K11").Value = 0.2150170648464164
$ws.Range("L11").Value = 0.515358361774744
$ws.Range("S11").Value = 0.03754266211604096
$ws.Range("G12").Value = 0.7532467532467533
$ws.Range("J12").Value = 0.2012987012987013
$ws.Range("K12").Value = 0.01298701298701299
$ws.Range("L12").Value = 0.01298701298701299
$ws.Range("S12").Value = 0.01948051948051948
$ws.Range("G13").Value = 0.6590909090909091
$ws.Range("J13").Value = 0.3409090909090909
$ws.Range("F15").Value = 0.03524229074889868
$ws.Range("H15").Value = 0.1365638766519824
$ws.Range("I15").Value = 0.08370044052863436
$ws.Range("J15").Value = 0.3612334801762114
$ws.Range("K15").Value = 0.07929515418502203
$ws.Range("M15").Value = 0.00881057268722467
$ws.Range("O15").Value = 0.1013215859030837
$ws.Range("S15").Value = 0.1938325991189427
$ws.Range("F16").Value = 0.02116402116402116
$ws.Range("H16").Value = 0.164021164021164
$ws.Range("I16").Value = 0.07936507936507936
$ws.Range("J16").Value = 0.4444444444444444
$ws.Range("K16").Value = 0.08465608465608465
$ws.Range("M16").Value = 0.02645502645502645
$ws.Range("O16").Value = 0.05291005291005291
$ws.Range("S16").Value = 0.126984126984127
$ws.Range("F17").Value = 0.01805869074492099
$ws.Range("H17").Value = 0.1467268623024831
$ws.Range("I17").Value = 0.08126410835214447
$ws.Range("J17").Value = 0.4379232505643341
$ws.Range("K17").Value = 0.1106094808126411
$ws.Range("M17").Value = 0.02257336343115124
$ws.Range("O17").Value = 0.0744920993227991
$ws.Range("S17").Value = 0.108352144469526
$ws.Range("F18").Value = 0.01941747572815534
$ws.Range("H18").Value = 0.1699029126213592
$ws.Range("I18").Value = 0.1067961165048544
$ws.Range("J18").Value = 0.412621359223301
$ws.Range("K18").Value = 0.1262135922330097
$ws.Range("M18").Value = 0.01456310679611651
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.06310679611650485
$ws.Range("S18").Value = 0.0825242718446602
$ws.Range("F19").Value = 0.01438848920863309
$ws.Range("H19").Value = 0.210431654676259
$ws.Range("I19").Value = 0.08633093525179857
$ws.Range("J19").Value = 0.381294964028777
$ws.Range("K19").Value = 0.1043165467625899
$ws.Range("M19").Value = 0.02338129496402878
$ws.Range("N19").Value = 0.0008992805755395684
$ws.Range("O19").Value = 0.07733812949640288
$ws.Range("S19").Value = 0.1016187050359712
